# Remove the stray "You might also like " boilerplate text that was
# scraped into the lyrics column (column C) for each song row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)   # Column C = Lyrics
    $val = $cell.Value2
    if ($val -ne $null -and $val.Contains("You might also like")) {
        $newVal = $val.Replace("You might also like ", "")
        $cell.Value2 = $newVal
    }
}
